$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instruments")

# Header E1: skos:description -> skos:definition
$ws.Range("E1").Value = "skos:definition"

# Fill column F (owl:sameAs target URLs) for each instrument row
$ws.Range("F2").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/CBQ_Family_J.pdf"
$ws.Range("F3").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DMQ_Family_J.pdf"
$ws.Range("F4").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/FSQ_Family_J.pdf"
$ws.Range("F5").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/HOQ_Family_J.pdf"
$ws.Range("F6").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/INQ_Family_J.pdf"
$ws.Range("F7").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/OCQ_Family_J.pdf"
$ws.Range("F8").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/SMQ_Family_J.pdf"
$ws.Range("F9").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/ACQ_J.pdf"
$ws.Range("F10").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/AUQ_J.pdf"
$ws.Range("F11").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/BPQ_J.pdf"
$ws.Range("F12").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/CDQ_J.pdf"
$ws.Range("F13").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DMQ_J.pdf"
$ws.Range("F14").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DEQ_J.pdf"
$ws.Range("F15").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DIQ_J.pdf"
$ws.Range("F16").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DBQ_J.pdf"
$ws.Range("F17").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DSQ_J.pdf"
$ws.Range("F18").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DLQ_J.pdf"
$ws.Range("F19").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/ECQ_J.pdf"
$ws.Range("F20").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/HIQ_J.pdf"
$ws.Range("F21").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/HEQ_J.pdf"
$ws.Range("F22").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/HUQ_J.pdf"
$ws.Range("F23").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/IMQ_J.pdf"
$ws.Range("F24").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/KIQ_J.pdf"
$ws.Range("F25").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/MCQ_J.pdf"
$ws.Range("F26").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/OCQ_J.pdf"
$ws.Range("F27").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/OHQ_J.pdf"
$ws.Range("F28").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/OSQ_J.pdf"
$ws.Range("F29").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/PAQ_J.pdf"
$ws.Range("F30").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/PFQ_J.pdf"
$ws.Range("F31").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/SLQ_J.pdf"
$ws.Range("F32").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/SMQ_J.pdf"
$ws.Range("F33").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/WHQ_J.pdf"
$ws.Range("F34").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/ALQ_ACASI_J.pdf"
$ws.Range("F35").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DUQ_ACASI_J.pdf"
$ws.Range("F36").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/SXQ_ACASI_J.pdf"
$ws.Range("F37").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/SMQ_ACASI_J.pdf"
$ws.Range("F38").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/ALQ_CAPI_J.pdf"
$ws.Range("F39").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/HSQ_CAPI_J.pdf"
$ws.Range("F40").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/DPQ_CAPI_J.pdf"
$ws.Range("F41").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/KIQ_CAPI_J.pdf"
$ws.Range("F42").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/PUQ_CAPI_J.pdf"
$ws.Range("F43").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/PAQ_CAPI_J.pdf"
$ws.Range("F44").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/RHQ_CAPI_J.pdf"
$ws.Range("F45").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/SXQ_CAPI_J.pdf"
$ws.Range("F46").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/SMQ_CAPI_J.pdf"
$ws.Range("F47").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/VTQ_CAPI_J.pdf"
$ws.Range("F48").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/WHQ_CAPI_J.pdf"
$ws.Range("F49").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/FCBS_PFU_Handcards_J.pdf"
$ws.Range("F50").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/FCBS-CBQ-J-508.pdf"
$ws.Range("F51").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2015-2016/questionnaires/DSA_I.pdf"
$ws.Range("F52").Value = "https://wwwn.cdc.gov/nchs/data/nhanes/2017-2018/questionnaires/Dietary-Post-Recall-2017-2018-508.pdf"
